$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Swap the store names for rows 3 and 4 (Ponta Negra <-> Vieiralves)
$ws.Range("A3").Value = "Bibi Cell Vieiralves"
$ws.Range("A4").Value = "Bibi Cell Ponta Negra"

# Row 2 - Bibi Cell Mundi
$ws.Range("C2").Value = 11735.4
$ws.Range("AG2").Value = 18540.55

# Row 3
$ws.Range("B3").Value = 3638
$ws.Range("C3").Value = 3280.25
$ws.Range("AG3").Value = 6918.25

# Row 4
$ws.Range("B4").Value = 4535.01
$ws.Range("C4").Value = 2016
$ws.Range("AG4").Value = 6551.01

# Row 5 - Bibi Cell Manauara
$ws.Range("B5").Value = 2756
$ws.Range("C5").Value = 3433
$ws.Range("AG5").Value = 6189

# Row 6 - total
$ws.Range("B6").Value = 17734.16
$ws.Range("C6").Value = 20464.65
$ws.Range("AG6").Value = 38198.81
